$p = $ppt.ActivePresentation

# Slide 2 ("What is Big Data?"): merge the "methods" + "." runs into
# a single run "methods." (retype the trailing substring so PowerPoint
# coalesces it with matching run formatting).
$s1 = $p.Slides.Item(2)
$shp1 = $s1.Shapes.Item(2)
$tr1 = $shp1.TextFrame.TextRange
$full1 = $tr1.Text
$idx1 = $full1.IndexOf("methods.")
$sub1 = $tr1.Characters($idx1 + 1, 8)
$sub1.Text = "methods."

# Slide 37 ("Some Drawbacks of Sharding"): merge the " is " + "fairly
# common, " runs into a single run " is fairly common, ".
$s2 = $p.Slides.Item(37)
$shp2 = $s2.Shapes.Item(2)
$tr2 = $shp2.TextFrame.TextRange
$full2 = $tr2.Text
$idx2 = $full2.IndexOf(" is fairly common, ")
$sub2 = $tr2.Characters($idx2 + 1, 19)
$sub2.Text = " is fairly common, "
